# Add a "rare_events" flag column (E) and a "report_comment" column (F)
# to the report_config example data, per commit:
#   "Update internal example data to include rare_events flag"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: rare_event_chart flag (Y/N) -----------------------------
$ws.Range("E1").Value = "rare_event_chart"
$ws.Range("E2").Value = "N"
$ws.Range("E3").Value = "N"
$ws.Range("E4").Value = "N"
$ws.Range("E5").Value = "N"
$ws.Range("E6").Value = "N"
$ws.Range("E7").Value = "N"
$ws.Range("E8").Value = "N"
$ws.Range("E9").Value = "N"
$ws.Range("E10").Value = "Y"

# --- New column F: report_comment (only populated on a couple of rows) -----
$ws.Range("F1").Value = "report_comment"
$ws.Range("F2").Value = "This is a comment about the attendances metric, which has been re-based as a demonstration.  This text is added via 'report_config.xlsx'"
$ws.Range("F7").Value = "Recent points demonstrate special-cause improvement.  Congratulations and carry on!"

# Match the saved selection state from the authored workbook.
[void]$ws.Range("G7").Select()
